$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Title "A" + " " + "slide" -> single run "A slide"
$tr1 = $s.Shapes.Item(1).TextFrame.TextRange
$tr1.Delete()
$tr1.Text = "A slide"

# Table cell "a" + " " + "table" -> single run "a table"
$tr2 = $s.Shapes.Item(3).Table.Cell(1, 2).Shape.TextFrame.TextRange
$tr2.Delete()
$tr2.Text = "a table"

# TextBox "Plus" + " " + "an" + " " + "image" -> single run "Plus an image"
$tr3 = $s.Shapes.Item(7).TextFrame.TextRange
$tr3.Delete()
$tr3.Text = "Plus an image"
